# Added and styled merged cells
# (Per the source diff, Sheet2 gets a new column F populated with the value 1
#  in rows 1-3; the active sheet/selection ends up on F3.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

$ws.Range("F1").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1

$ws.Range("F3").Select()
